$wb = $excel.ActiveWorkbook

$sheetNames = @("2025", "2030", "2035", "2040", "2045", "2050")

$values2025 = @(0, 290.0628494009878, 0, 0, 29049.07128553874, 0, 8095.925712662051, 0, 14945.10834652955, 0, 0, 50998.86069102, 11228.70813999, 7234.0658054822, 6703.624349245061)
$values2030 = @(219.6160489230463, 3803.736742006062, 0, 0, 45497.55843345862, 0, 8095.925712662051, 0, 31114.43531462794, 0, 0, 60434.98124678315, 17372.009741075, 9195.867044489814, 7876.358949184817)
$values2035 = @(2152.642131479708, 5727.815297363306, 0, 0, 57498.34554108262, 0, 8095.925712662051, 0, 49824.6994397517, 0, 0, 60434.98124678315, 23386.44280563801, 13671.2814398324, 12393.80935914757)
$values2040 = @(2152.642131479708, 5727.815297363306, 0, 0, 57498.34554108262, 0, 8095.925712662051, 0, 49824.6994397517, 0, 0, 60434.98124678315, 23386.44280563801, 13671.2814398324, 12393.80935914757)
$values2045 = @(2152.642131479708, 5727.815297363306, 0, 0, 57498.34554108262, 0, 8095.925712662051, 0, 49824.6994397517, 0, 0, 60434.98124678315, 23386.44280563801, 13671.2814398324, 14414.6765450033)
$values2050 = @(2152.642131479708, 5727.815297363306, 0, 0, 57498.34554108262, 0, 8095.925712662051, 0, 49824.6994397517, 0, 0, 60434.98124678315, 23386.44280563801, 13671.2814398324, 14414.6765450033)

$allValues = @($values2025, $values2030, $values2035, $values2040, $values2045, $values2050)

for ($s = 0; $s -lt $sheetNames.Count; $s++) {
    $ws = $wb.Worksheets($sheetNames[$s])
    $values = $allValues[$s]
    for ($i = 0; $i -lt $values.Count; $i++) {
        $col = $i + 1
        $ws.Cells.Item(2, $col).Value = $values[$i]
    }
}
